# Weekly driver report update for 2025-05-05
# Updates the "Good Drivers" table (rows 12-33) on the "Driver Summary"
# sheet with refreshed counts/dates for this week's pull. Date-looking
# strings in column J are written with a leading apostrophe so Excel
# keeps them as literal text instead of auto-converting to date serials.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.250.0.4"
$ws.Range("B12").Value = 1869842
$ws.Range("C12").Value = 3420
$ws.Range("D12").Value = 4386
$ws.Range("E12").Value = 2436
$ws.Range("F12").Value = 1877648
$ws.Range("H12").Value = "22.250.0.4"
$ws.Range("J12").Value = "'2023-07-25"
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.40.0.7"
$ws.Range("B13").Value = 8170878
$ws.Range("C13").Value = 15867
$ws.Range("D13").Value = 13188
$ws.Range("E13").Value = 18436
$ws.Range("F13").Value = 8199933
$ws.Range("H13").Value = "22.40.0.7"
$ws.Range("J13").Value = "'2021-09-18"
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.190.0.4"
$ws.Range("B14").Value = 1611822
$ws.Range("C14").Value = 4793
$ws.Range("D14").Value = 1136
$ws.Range("E14").Value = 2996
$ws.Range("F14").Value = 1617751
$ws.Range("H14").Value = "22.190.0.4"
$ws.Range("J14").Value = "'2022-11-22"
$ws.Range("B15").Value = 298304
$ws.Range("C15").Value = 925
$ws.Range("E15").Value = 420
$ws.Range("F15").Value = 299532
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.10.0.8"
$ws.Range("B16").Value = 381616
$ws.Range("C16").Value = 1106
$ws.Range("D16").Value = 582
$ws.Range("E16").Value = 601
$ws.Range("F16").Value = 383304
$ws.Range("H16").Value = "23.10.0.8"
$ws.Range("I16").Value = 99.6
$ws.Range("J16").Value = "'2023-10-30"
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.90.0.5"
$ws.Range("B17").Value = 335610
$ws.Range("C17").Value = 495
$ws.Range("D17").Value = 530
$ws.Range("E17").Value = 414
$ws.Range("F17").Value = 336635
$ws.Range("H17").Value = "22.90.0.5"
$ws.Range("J17").Value = "'2021-09-26"
$ws.Range("A18").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.220.0.4"
$ws.Range("B18").Value = 750778
$ws.Range("C18").Value = 1655
$ws.Range("D18").Value = 507
$ws.Range("E18").Value = 1290
$ws.Range("F18").Value = 752940
$ws.Range("H18").Value = "22.220.0.4"
$ws.Range("J18").Value = "'2023-03-28"
$ws.Range("A19").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.20.0.6"
$ws.Range("B19").Value = 1021828
$ws.Range("C19").Value = 2379
$ws.Range("D19").Value = 1006
$ws.Range("E19").Value = 1021
$ws.Range("F19").Value = 1025213
$ws.Range("H19").Value = "22.20.0.6"
$ws.Range("J19").Value = "'2020-11-29"
$ws.Range("A20").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.60.0.6"
$ws.Range("B20").Value = 100380
$ws.Range("C20").Value = 264
$ws.Range("D20").Value = 30
$ws.Range("E20").Value = 267
$ws.Range("F20").Value = 100674
$ws.Range("H20").Value = "22.60.0.6"
$ws.Range("J20").Value = "'2021-05-26"
$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.80.1.3"
$ws.Range("B21").Value = 195464
$ws.Range("C21").Value = 430
$ws.Range("D21").Value = 70
$ws.Range("E21").Value = 401
$ws.Range("F21").Value = 195964
$ws.Range("H21").Value = "23.80.1.3"
$ws.Range("J21").Value = "'2024-09-03"
$ws.Range("A22").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.0.0.6"
$ws.Range("B22").Value = 1833336
$ws.Range("C22").Value = 2376
$ws.Range("D22").Value = 656
$ws.Range("E22").Value = 1899
$ws.Range("F22").Value = 1836368
$ws.Range("H22").Value = "22.0.0.6"
$ws.Range("J22").Value = "'2020-09-16"
$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.170.0.3"
$ws.Range("B23").Value = 1237156
$ws.Range("C23").Value = 2499
$ws.Range("D23").Value = 583
$ws.Range("E23").Value = 1951
$ws.Range("F23").Value = 1240238
$ws.Range("H23").Value = "22.170.0.3"
$ws.Range("J23").Value = "'2022-08-28"
$ws.Range("A24").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.80.0.4"
$ws.Range("B24").Value = 63994
$ws.Range("C24").Value = 66
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = 68
$ws.Range("F24").Value = 64101
$ws.Range("H24").Value = "21.80.0.4"
$ws.Range("J24").Value = "'2020-01-29"
$ws.Range("A25").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B25").Value = 81417
$ws.Range("C25").Value = 107
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = 158
$ws.Range("F25").Value = 81549
$ws.Range("H25").Value = "22.80.0.9"
$ws.Range("J25").Value = "'2021-08-18"
$ws.Range("A26").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.0.1.1"
$ws.Range("B26").Value = 1400496
$ws.Range("C26").Value = 2029
$ws.Range("D26").Value = 1075
$ws.Range("E26").Value = 3210
$ws.Range("F26").Value = 1403600
$ws.Range("H26").Value = "22.0.1.1"
$ws.Range("J26").Value = "'2020-09-28"
$ws.Range("A27").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.80.2.1"
$ws.Range("B27").Value = 393045
$ws.Range("C27").Value = 453
$ws.Range("D27").Value = 216
$ws.Range("E27").Value = 1251
$ws.Range("F27").Value = 393714
$ws.Range("H27").Value = "21.80.2.1"
$ws.Range("I27").Value = 99.8
$ws.Range("J27").Value = "'2020-02-24"
$ws.Range("B28").Value = 509717
$ws.Range("C28").Value = 522
$ws.Range("E28").Value = 590
$ws.Range("F28").Value = 510417
$ws.Range("A29").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.40.1.3"
$ws.Range("B29").Value = 11362
$ws.Range("C29").Value = 7
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 30
$ws.Range("F29").Value = 11369
$ws.Range("H29").Value = "21.40.1.3"
$ws.Range("I29").Value = 99.9
$ws.Range("J29").Value = ""
$ws.Range("A30").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B30").Value = 75454
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 79
$ws.Range("F30").Value = 75457
$ws.Range("H30").Value = "21.110.3.2"
$ws.Range("J30").Value = "'2020-08-05"
$ws.Range("B32").Value = 121310
$ws.Range("C32").Value = 26
$ws.Range("E32").Value = 154
$ws.Range("F32").Value = 121350
$ws.Range("B33").Value = 35363
$ws.Range("E33").Value = 81
$ws.Range("F33").Value = 35377
